$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write the final table contents for rows 2-7, columns A-T.
# (Sending cluster / Ligand symbol / Receptor symbol / Target cluster / metrics)

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Il18"
$ws.Cells.Item(2, 3).Value = "Il18rap"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 2.878986666666667
$ws.Cells.Item(2, 8).Value = 8.63696
$ws.Cells.Item(2, 9).Value = 0.2331567682967092
$ws.Cells.Item(2, 10).Value = 0.2331567682967092
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 25.939101
$ws.Cells.Item(2, 14).Value = 77.81730300000001
$ws.Cells.Item(2, 15).Value = 0.9986086494543361
$ws.Cells.Item(2, 16).Value = 0.9986086494543361
$ws.Cells.Item(2, 17).Value = 74.67832592432
$ws.Cells.Item(2, 18).Value = 672.1049333188801
$ws.Cells.Item(2, 19).Value = 0.2328323654999143
$ws.Cells.Item(2, 20).Value = 0.2328323654999143

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Il18"
$ws.Cells.Item(3, 3).Value = "Il18rap"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 2.878986666666667
$ws.Cells.Item(3, 8).Value = 8.63696
$ws.Cells.Item(3, 9).Value = 0.2331567682967092
$ws.Cells.Item(3, 10).Value = 0.2331567682967092
$ws.Cells.Item(3, 11).Value = 1
$ws.Cells.Item(3, 12).Value = 0.3333333333333333
$ws.Cells.Item(3, 13).Value = 0.03614066666666667
$ws.Cells.Item(3, 14).Value = 0.108422
$ws.Cells.Item(3, 15).Value = 0.001391350545663835
$ws.Cells.Item(3, 16).Value = 0.001391350545663835
$ws.Cells.Item(3, 17).Value = 0.1040484974577778
$ws.Cells.Item(3, 18).Value = 0.9364364771200001
$ws.Cells.Item(3, 19).Value = 0.0003244027967948428
$ws.Cells.Item(3, 20).Value = 0.0003244027967948428

# Row 4
$ws.Cells.Item(4, 1).Value = "FAPs"
$ws.Cells.Item(4, 2).Value = "Il18"
$ws.Cells.Item(4, 3).Value = "Il18rap"
$ws.Cells.Item(4, 4).Value = "ECs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 8.325812333333333
$ws.Cells.Item(4, 8).Value = 24.977437
$ws.Cells.Item(4, 9).Value = 0.6742717913773656
$ws.Cells.Item(4, 10).Value = 0.6742717913773655
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 25.939101
$ws.Cells.Item(4, 14).Value = 77.81730300000001
$ws.Cells.Item(4, 15).Value = 0.9986086494543361
$ws.Cells.Item(4, 16).Value = 0.9986086494543361
$ws.Cells.Item(4, 17).Value = 215.964087021379
$ws.Cells.Item(4, 18).Value = 1943.676783192411
$ws.Cells.Item(4, 19).Value = 0.673333642952507
$ws.Cells.Item(4, 20).Value = 0.6733336429525069

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Il18"
$ws.Cells.Item(5, 3).Value = "Il18rap"
$ws.Cells.Item(5, 4).Value = "FAPs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 8.325812333333333
$ws.Cells.Item(5, 8).Value = 24.977437
$ws.Cells.Item(5, 9).Value = 0.6742717913773656
$ws.Cells.Item(5, 10).Value = 0.6742717913773655
$ws.Cells.Item(5, 11).Value = 1
$ws.Cells.Item(5, 12).Value = 0.3333333333333333
$ws.Cells.Item(5, 13).Value = 0.03614066666666667
$ws.Cells.Item(5, 14).Value = 0.108422
$ws.Cells.Item(5, 15).Value = 0.001391350545663835
$ws.Cells.Item(5, 16).Value = 0.001391350545663835
$ws.Cells.Item(5, 17).Value = 0.3009004082682222
$ws.Cells.Item(5, 18).Value = 2.708103674414
$ws.Cells.Item(5, 19).Value = 0.0009381484248586295
$ws.Cells.Item(5, 20).Value = 0.0009381484248586293

# Row 6
$ws.Cells.Item(6, 1).Value = "sCs"
$ws.Cells.Item(6, 2).Value = "Il18"
$ws.Cells.Item(6, 3).Value = "Il18rap"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 1.143059
$ws.Cells.Item(6, 8).Value = 3.429177
$ws.Cells.Item(6, 9).Value = 0.09257144032592537
$ws.Cells.Item(6, 10).Value = 0.09257144032592536
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 25.939101
$ws.Cells.Item(6, 14).Value = 77.81730300000001
$ws.Cells.Item(6, 15).Value = 0.9986086494543361
$ws.Cells.Item(6, 16).Value = 0.9986086494543361
$ws.Cells.Item(6, 17).Value = 29.649922849959
$ws.Cells.Item(6, 18).Value = 266.849305649631
$ws.Cells.Item(6, 19).Value = 0.09244264100191502
$ws.Cells.Item(6, 20).Value = 0.092442641001915

# Row 7
$ws.Cells.Item(7, 1).Value = "sCs"
$ws.Cells.Item(7, 2).Value = "Il18"
$ws.Cells.Item(7, 3).Value = "Il18rap"
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 1.143059
$ws.Cells.Item(7, 8).Value = 3.429177
$ws.Cells.Item(7, 9).Value = 0.09257144032592537
$ws.Cells.Item(7, 10).Value = 0.09257144032592536
$ws.Cells.Item(7, 11).Value = 1
$ws.Cells.Item(7, 12).Value = 0.3333333333333333
$ws.Cells.Item(7, 13).Value = 0.03614066666666667
$ws.Cells.Item(7, 14).Value = 0.108422
$ws.Cells.Item(7, 15).Value = 0.001391350545663835
$ws.Cells.Item(7, 16).Value = 0.001391350545663835
$ws.Cells.Item(7, 17).Value = 0.04131091429933334
$ws.Cells.Item(7, 18).Value = 0.371798228694
$ws.Cells.Item(7, 19).Value = 0.0001287993240103635
$ws.Cells.Item(7, 20).Value = 0.0001287993240103634
